$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update name/position/time for Fazliddin Xamdamov (Assistant)
$ws.Range("A2").Value = "Fazliddin Xamdamov"
$ws.Range("B2").Value = "Assistant"
$ws.Range("C2").Value = "Osh"
$ws.Range("D2").Value = "2022-03-14T12:49:23.242951"

# Row 3: second food log entry for the same person, later timestamp
$ws.Range("A3").Value = "Fazliddin Xamdamov"
$ws.Range("B3").Value = "Assistant"
$ws.Range("C3").Value = "Osh"
$ws.Range("D3").Value = "2022-03-14T13:22:13.154746"

# Row 4 no longer exists - clear it out entirely
$ws.Range("A4:D4").Clear()
